$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("Reporting Managers") for the new
# "Reporting ManagerId" field. This shifts the existing C:L columns to D:M.
$ws.Columns.Item(3).Insert() | Out-Null

# New header cell for the inserted column.
$ws.Cells.Item(1, 3).Value = "Reporting ManagerId"

# Give the new column a sensible manual width (matches the "Emp Name" column
# width, but set explicitly rather than via best-fit/autofit).
$ws.Columns.Item(3).ColumnWidth = 9.67

# Selection collapses back down to just the single new header cell.
$ws.Range("C1").Select() | Out-Null
